{"js": "// Cover page edit: set the real dissertation title, subtitle and supervisor\n// name (replacing the placeholder \"TBD\" content), per:\n//   \"Insert initial information (title, supervisor, author)\"\n//\n// 1) Title paragraph:\n//      \"Multi-Platform large dense matrices multiplication for performance\n//       enhancing\"\n//    ->  three runs (same run formatting throughout):\n//      \"Design of a multi-sensor\" + \" \" + \"apparatus for  forestry robotics\"\n//\n// 2) Subtitle paragraph:\n//      \"Subt\u00edtulo\" (run #1, sz 40) + \" TBD\" (run #2, sz 32)\n//    -> single run (keeps the sz-32 run formatting):\n//      \"A case study for forest 3D mapping \"\n//\n// 3) Supervisor paragraph: swap the supervisor's name in place, keeping the\n//    rest of the sentence and the run formatting untouched.\n\n// ---- 1) Title ---------------------------------------------------------\nconst titleHits = context.document.body.search(\n  \"Multi-Platform large dense matrices multiplication for performance enhancing\",\n  { matchCase: true }\n);\ntitleHits.load(\"items\");\nawait context.sync();\n\nif (titleHits.items.length === 0) {\n  throw new Error(\"Title paragraph text not found\");\n}\nconst titlePara = titleHits.items[0].paragraphs.getFirst();\n\n// Run properties are identical to the ones on the original title run, so we\n// rebuild the paragraph via OOXML to get three discrete <w:r> elements\n// (each run keeps the paragraph's existing formatting).\nconst titleRunProps =\n  '<w:rFonts w:eastAsia=\"Times New Roman\" w:cs=\"Times New Roman\" w:ascii=\"Arrus BT Roman\" w:hAnsi=\"Arrus BT Roman\"/>' +\n  \"<w:b/>\" +\n  \"<w:smallCaps/>\" +\n  '<w:color w:val=\"6B583E\"/>' +\n  '<w:kern w:val=\"0\"/>' +\n  '<w:sz w:val=\"44\"/>' +\n  '<w:szCs w:val=\"44\"/>' +\n  '<w:lang w:val=\"pt-PT\" w:eastAsia=\"pt-PT\" w:bidi=\"ar-SA\"/>';\n\nconst titleOoxml =\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>' +\n  \"<w:p>\" +\n  '<w:pPr><w:pStyle w:val=\"Normal\"/><w:jc w:val=\"center\"/><w:rPr>' +\n  titleRunProps +\n  \"</w:rPr></w:pPr>\" +\n  \"<w:r><w:rPr>\" +\n  titleRunProps +\n  \"</w:rPr><w:t>Design of a multi-sensor</w:t></w:r>\" +\n  \"<w:r><w:rPr>\" +\n  titleRunProps +\n  '</w:rPr><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n  \"<w:r><w:rPr>\" +\n  titleRunProps +\n  \"</w:rPr><w:t>apparatus for  forestry robotics</w:t></w:r>\" +\n  \"</w:p>\" +\n  \"</w:body></w:document>\" +\n  \"</pkg:xmlData></pkg:part></pkg:package>\";\n\ntitlePara.insertOoxml(titleOoxml, \"Replace\");\nawait context.sync();\n\n// ---- 2) Subtitle --------------------------------------------------------\nconst subtitleHits = context.document.body.search(\"Subt\u00edtulo\", { matchCase: true });\nsubtitleHits.load(\"items\");\nawait context.sync();\n\nif (subtitleHits.items.length === 0) {\n  throw new Error(\"Subtitle paragraph text not found\");\n}\nconst subtitlePara = subtitleHits.items[0].paragraphs.getFirst();\n\n// Keep the paragraph mark / paragraph-level formatting and the formatting of\n// the second run (sz 32); the first run (\"Subt\u00edtulo\", sz 40) is dropped.\nconst subtitleParaRunProps =\n  '<w:rFonts w:ascii=\"Arrus BT Roman\" w:hAnsi=\"Arrus BT Roman\" w:eastAsia=\"Times New Roman\" w:cs=\"Times New Roman\"/>' +\n  '<w:color w:val=\"6B583E\"/>' +\n  '<w:sz w:val=\"32\"/>' +\n  '<w:szCs w:val=\"32\"/>' +\n  '<w:lang w:eastAsia=\"pt-PT\"/>';\n\nconst subtitleRunProps =\n  '<w:rFonts w:eastAsia=\"Times New Roman\" w:cs=\"Times New Roman\" w:ascii=\"Arrus BT Roman\" w:hAnsi=\"Arrus BT Roman\"/>' +\n  \"<w:smallCaps/>\" +\n  '<w:color w:val=\"6B583E\"/>' +\n  '<w:sz w:val=\"32\"/>' +\n  '<w:szCs w:val=\"32\"/>' +\n  '<w:lang w:eastAsia=\"pt-PT\"/>';\n\nconst subtitleOoxml =\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>' +\n  \"<w:p>\" +\n  '<w:pPr><w:pStyle w:val=\"Normal\"/><w:jc w:val=\"center\"/><w:rPr>' +\n  subtitleParaRunProps +\n  \"</w:rPr></w:pPr>\" +\n  \"<w:r><w:rPr>\" +\n  subtitleRunProps +\n  '</w:rPr><w:t xml:space=\"preserve\">A case study for forest 3D mapping </w:t></w:r>' +\n  \"</w:p>\" +\n  \"</w:body></w:document>\" +\n  \"</pkg:xmlData></pkg:part></pkg:package>\";\n\nsubtitlePara.insertOoxml(subtitleOoxml, \"Replace\");\nawait context.sync();\n\n// ---- 3) Supervisor name --------------------------------------------------\nconst supervisorHits = context.document.body.search(\"V\u00edtor Manuel Mendes da Silva\", {\n  matchCase: true,\n});\nsupervisorHits.load(\"items\");\nawait context.sync();\n\nif (supervisorHits.items.length === 0) {\n  throw new Error(\"Supervisor name not found\");\n}\nconst supervisorPara = supervisorHits.items[0].paragraphs.getFirst();\nsupervisorPara.load(\"text\");\nawait context.sync();\n\nconst updatedSentence = supervisorPara.text.replace(\n  \"V\u00edtor Manuel Mendes da Silva\",\n  \"David Bina Siassipour Portugal\"\n);\n// insertText(\"Replace\") on the whole paragraph range keeps the existing run\n// formatting (single run, unchanged rPr) while swapping the text.\nsupervisorPara.insertText(updatedSentence, \"Replace\");\nawait context.sync();\n", "ps1": "# Cover page edit: set the real dissertation title, subtitle and supervisor\n# name (replacing the placeholder \"TBD\" content), per:\n#   \"Insert initial information (title, supervisor, author)\"\n#\n# 1) Title paragraph:\n#      \"Multi-Platform large dense matrices multiplication for performance\n#       enhancing\"\n#    ->  three runs (same run formatting throughout):\n#      \"Design of a multi-sensor\" + \" \" + \"apparatus for  forestry robotics\"\n#\n# 2) Subtitle paragraph:\n#      \"Subt\u00edtulo\" (run #1, sz 40) + \" TBD\" (run #2, sz 32)\n#    -> single run (keeps the sz-32 run formatting):\n#      \"A case study for forest 3D mapping \"\n#\n# 3) Supervisor paragraph: swap the supervisor's name in place, keeping the\n#    rest of the sentence and the run formatting untouched.\n\n$d = $word.ActiveDocument\n\nfunction Get-ParagraphRangeByText($searchText) {\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Text = $searchText\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $found = $find.Execute()\n    if (-not $found) {\n        throw \"Could not find text: $searchText\"\n    }\n    # Grow the found hit out to the paragraph that contains it.\n    $rng.Expand(4) | Out-Null  # wdParagraph\n    return $rng\n}\n\n# ---- 1) Title -------------------------------------------------------------\n$titleRange = Get-ParagraphRangeByText(\"Multi-Platform large dense matrices multiplication for performance enhancing\")\n\n$titleRunProps = '<w:rFonts w:eastAsia=\"Times New Roman\" w:cs=\"Times New Roman\" w:ascii=\"Arrus BT Roman\" w:hAnsi=\"Arrus BT Roman\"/>' + `\n    '<w:b/>' + `\n    '<w:smallCaps/>' + `\n    '<w:color w:val=\"6B583E\"/>' + `\n    '<w:kern w:val=\"0\"/>' + `\n    '<w:sz w:val=\"44\"/>' + `\n    '<w:szCs w:val=\"44\"/>' + `\n    '<w:lang w:val=\"pt-PT\" w:eastAsia=\"pt-PT\" w:bidi=\"ar-SA\"/>'\n\n$titleOoxml = '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' + `\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData>' + `\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>' + `\n    '<w:p>' + `\n    '<w:pPr><w:pStyle w:val=\"Normal\"/><w:jc w:val=\"center\"/><w:rPr>' + $titleRunProps + '</w:rPr></w:pPr>' + `\n    '<w:r><w:rPr>' + $titleRunProps + '</w:rPr><w:t>Design of a multi-sensor</w:t></w:r>' + `\n    '<w:r><w:rPr>' + $titleRunProps + '</w:rPr><w:t xml:space=\"preserve\"> </w:t></w:r>' + `\n    '<w:r><w:rPr>' + $titleRunProps + '</w:rPr><w:t>apparatus for  forestry robotics</w:t></w:r>' + `\n    '</w:p>' + `\n    '</w:body></w:document>' + `\n    '</pkg:xmlData></pkg:part></pkg:package>'\n\n$titleRange.InsertXML($titleOoxml)\n\n# ---- 2) Subtitle ------------------------------------------------------------\n$subtitleRange = Get-ParagraphRangeByText(\"Subt\u00edtulo\")\n\n$subtitleParaRunProps = '<w:rFonts w:ascii=\"Arrus BT Roman\" w:hAnsi=\"Arrus BT Roman\" w:eastAsia=\"Times New Roman\" w:cs=\"Times New Roman\"/>' + `\n    '<w:color w:val=\"6B583E\"/>' + `\n    '<w:sz w:val=\"32\"/>' + `\n    '<w:szCs w:val=\"32\"/>' + `\n    '<w:lang w:eastAsia=\"pt-PT\"/>'\n\n$subtitleRunProps = '<w:rFonts w:eastAsia=\"Times New Roman\" w:cs=\"Times New Roman\" w:ascii=\"Arrus BT Roman\" w:hAnsi=\"Arrus BT Roman\"/>' + `\n    '<w:smallCaps/>' + `\n    '<w:color w:val=\"6B583E\"/>' + `\n    '<w:sz w:val=\"32\"/>' + `\n    '<w:szCs w:val=\"32\"/>' + `\n    '<w:lang w:eastAsia=\"pt-PT\"/>'\n\n$subtitleOoxml = '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' + `\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData>' + `\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>' + `\n    '<w:p>' + `\n    '<w:pPr><w:pStyle w:val=\"Normal\"/><w:jc w:val=\"center\"/><w:rPr>' + $subtitleParaRunProps + '</w:rPr></w:pPr>' + `\n    '<w:r><w:rPr>' + $subtitleRunProps + '</w:rPr><w:t xml:space=\"preserve\">A case study for forest 3D mapping </w:t></w:r>' + `\n    '</w:p>' + `\n    '</w:body></w:document>' + `\n    '</pkg:xmlData></pkg:part></pkg:package>'\n\n$subtitleRange.InsertXML($subtitleOoxml)\n\n# ---- 3) Supervisor name -----------------------------------------------------\n$rng = $d.Content\n$find = $rng.Find\n$find.ClearFormatting()\n$find.Text = \"V\u00edtor Manuel Mendes da Silva\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"David Bina Siassipour Portugal\"\n$find.MatchCase = $true\n# wdFindContinue = 1, wdReplaceOne = 1\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 1) | Out-Null\n"}
